$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.561.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.684.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "664.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.74%  "
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.696.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.555.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.117"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "16.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "467.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.645"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "79.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.830.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  +2.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.59%  "
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("E34").Value = "  +3.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.676.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.45%  "
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "178.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.76%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0900"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.64%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.50%  "
$ws.Range("E48").Value = "  -3.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("E51").Value = "  -1.28%  "
